$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = '''27.405.36'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '''  -2.60%  '
$ws.Range("E2").Style = "Normal"
# Row 3
$ws.Range("D3").Value = '''1.738.96'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '''  -3.44%  '
$ws.Range("E3").Style = "Normal"
# Row 4
$ws.Range("D4").Value = '''1.003'
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = '''  -0.21%  '
$ws.Range("E4").Style = "Normal"
# Row 5
$ws.Range("D5").Value = '''324.08'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '''  -4.28%  '
$ws.Range("E5").Style = "Normal"
# Row 6
$ws.Range("D6").Value = '''1.002'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '''  +0.17%  '
$ws.Range("E6").Style = "Normal"
# Row 7
$ws.Range("D7").Value = '''0.4266'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '''  -12.82%  '
$ws.Range("E7").Style = "Normal"
# Row 8
$ws.Range("D8").Value = '''0.3629'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '''  -2.65%  '
$ws.Range("E8").Style = "Normal"
# Row 9
$ws.Range("D9").Value = '''44.96'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '''  -1.41%  '
$ws.Range("E9").Style = "Normal"
# Row 10
$ws.Range("D10").Value = '''1.120'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '''  -2.18%  '
$ws.Range("E10").Style = "Normal"
# Row 11
$ws.Range("D11").Value = '''0.07381'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '''  -4.33%  '
$ws.Range("E11").Style = "Normal"
# Row 12
$ws.Range("D12").Value = '''1.001'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '''  -0.33%  '
$ws.Range("E12").Style = "Normal"
# Row 13
$ws.Range("D13").Value = '''21.66'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '''  -3.96%  '
$ws.Range("E13").Style = "Normal"
# Row 14
$ws.Range("D14").Value = '''6.069'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '''  -4.05%  '
$ws.Range("E14").Style = "Normal"
# Row 15
$ws.Range("D15").Value = '''7.185'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '''  -1.82%  '
$ws.Range("E15").Style = "Normal"
# Row 16
$ws.Range("D16").Value = '''1.735.02'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '''  -3.46%  '
$ws.Range("E16").Style = "Normal"
# Row 17
$ws.Range("D17").Value = '''0.00001059'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '''  -3.53%  '
$ws.Range("E17").Style = "Normal"
# Row 18
$ws.Range("D18").Value = '''85.23'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '''  +3.63%  '
$ws.Range("E18").Style = "Normal"
# Row 19
$ws.Range("D19").Value = '''0.05938'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '''  -11.79%  '
$ws.Range("E19").Style = "Normal"
# Row 20
$ws.Range("D20").Value = '''1.002'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '''  +0.18%  '
$ws.Range("E20").Style = "Normal"
# Row 21
$ws.Range("D21").Value = '''16.81'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '''  -3.38%  '
$ws.Range("E21").Style = "Normal"
# Row 22
$ws.Range("B22").Value = '''BitDAO'
$ws.Range("B22").Style = "Normal"
$ws.Range("C22").Value = '''https://coinranking.com/coin/N2IgQ9Xme+bitdao-bit'
$ws.Range("C22").Style = "Normal"
$ws.Range("D22").Value = '''0.5279'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '''  -4.51%  '
$ws.Range("E22").Style = "Normal"
# Row 23
$ws.Range("B23").Value = '''Uniswap'
$ws.Range("B23").Style = "Normal"
$ws.Range("C23").Value = '''https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
$ws.Range("C23").Style = "Normal"
$ws.Range("D23").Value = '''6.026'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '''  -6.28%  '
$ws.Range("E23").Style = "Normal"
# Row 24
$ws.Range("B24").Value = '''WrappedBTC'
$ws.Range("B24").Style = "Normal"
$ws.Range("C24").Value = '''https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'
$ws.Range("C24").Style = "Normal"
$ws.Range("D24").Value = '''27.409.76'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '''  -2.70%  '
$ws.Range("E24").Style = "Normal"
# Row 25
$ws.Range("B25").Value = '''Cosmos'
$ws.Range("B25").Style = "Normal"
$ws.Range("C25").Value = '''https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range("C25").Style = "Normal"
$ws.Range("D25").Value = '''11.27'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '''  -6.04%  '
$ws.Range("E25").Style = "Normal"
# Row 26
$ws.Range("B26").Value = '''Toncoin'
$ws.Range("B26").Style = "Normal"
$ws.Range("C26").Value = '''https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range("C26").Style = "Normal"
$ws.Range("D26").Value = '''2.399'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '''  -0.30%  '
$ws.Range("E26").Style = "Normal"
# Row 27
$ws.Range("B27").Value = '''EthereumClassic'
$ws.Range("B27").Style = "Normal"
$ws.Range("C27").Value = '''https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range("C27").Style = "Normal"
$ws.Range("D27").Value = '''19.95'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '''  -4.78%  '
$ws.Range("E27").Style = "Normal"
# Row 28
$ws.Range("B28").Value = '''LidoDAOToken'
$ws.Range("B28").Style = "Normal"
$ws.Range("C28").Value = '''https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range("C28").Style = "Normal"
$ws.Range("D28").Value = '''2.348'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '''  -2.69%  '
$ws.Range("E28").Style = "Normal"
# Row 29
$ws.Range("B29").Value = '''Monero'
$ws.Range("B29").Style = "Normal"
$ws.Range("C29").Value = '''https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range("C29").Style = "Normal"
$ws.Range("D29").Value = '''148.70'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '''  -1.80%  '
$ws.Range("E29").Style = "Normal"
# Row 30
$ws.Range("B30").Value = '''WrappedliquidstakedEther2.0'
$ws.Range("B30").Style = "Normal"
$ws.Range("C30").Value = '''https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range("C30").Style = "Normal"
$ws.Range("D30").Value = '''1.934.09'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '''  -3.56%  '
$ws.Range("E30").Style = "Normal"
# Row 31
$ws.Range("B31").Value = '''ImmutableX'
$ws.Range("B31").Style = "Normal"
$ws.Range("C31").Value = '''https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range("C31").Style = "Normal"
$ws.Range("D31").Value = '''1.247'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '''  -2.18%  '
$ws.Range("E31").Style = "Normal"
# Row 32
$ws.Range("B32").Value = '''BitcoinCash'
$ws.Range("B32").Style = "Normal"
$ws.Range("C32").Value = '''https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$ws.Range("C32").Style = "Normal"
$ws.Range("D32").Value = '''125.93'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '''  -6.25%  '
$ws.Range("E32").Style = "Normal"
# Row 33
$ws.Range("B33").Value = '''HuobiToken'
$ws.Range("B33").Style = "Normal"
$ws.Range("C33").Value = '''https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$ws.Range("C33").Style = "Normal"
$ws.Range("D33").Value = '''3.734'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '''  -7.72%  '
$ws.Range("E33").Style = "Normal"
# Row 34
$ws.Range("D34").Value = '''5.580'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '''  -6.15%  '
$ws.Range("E34").Style = "Normal"
# Row 35
$ws.Range("B35").Value = '''Stellar'
$ws.Range("B35").Style = "Normal"
$ws.Range("C35").Value = '''https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range("C35").Style = "Normal"
$ws.Range("D35").Value = '''0.09011'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '''  -8.89%  '
$ws.Range("E35").Style = "Normal"
# Row 36
$ws.Range("B36").Value = '''Aptos'
$ws.Range("B36").Style = "Normal"
$ws.Range("C36").Value = '''https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range("C36").Style = "Normal"
$ws.Range("D36").Value = '''12.37'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '''  +1.22%  '
$ws.Range("E36").Style = "Normal"
# Row 37
$ws.Range("B37").Value = '''Algorand'
$ws.Range("B37").Style = "Normal"
$ws.Range("C37").Value = '''https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range("C37").Style = "Normal"
$ws.Range("D37").Value = '''0.2165'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '''  -2.65%  '
$ws.Range("E37").Style = "Normal"
# Row 38
$ws.Range("B38").Value = '''VeChain'
$ws.Range("B38").Style = "Normal"
$ws.Range("C38").Value = '''https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("C38").Style = "Normal"
$ws.Range("D38").Value = '''0.02277'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '''  -4.48%  '
$ws.Range("E38").Style = "Normal"
# Row 39
$ws.Range("B39").Value = '''Hedera'
$ws.Range("B39").Style = "Normal"
$ws.Range("C39").Value = '''https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range("C39").Style = "Normal"
$ws.Range("D39").Value = '''0.06118'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '''  -4.32%  '
$ws.Range("E39").Style = "Normal"
# Row 40
$ws.Range("B40").Value = '''TheSandbox'
$ws.Range("B40").Style = "Normal"
$ws.Range("C40").Value = '''https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$ws.Range("C40").Style = "Normal"
$ws.Range("D40").Value = '''0.6434'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '''  -4.13%  '
$ws.Range("E40").Style = "Normal"
# Row 41
$ws.Range("B41").Value = '''InternetComputer(DFINITY)'
$ws.Range("B41").Style = "Normal"
$ws.Range("C41").Value = '''https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range("C41").Style = "Normal"
$ws.Range("D41").Value = '''4.993'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '''  -4.55%  '
$ws.Range("E41").Style = "Normal"
# Row 42
$ws.Range("B42").Value = '''TrustWalletToken'
$ws.Range("B42").Style = "Normal"
$ws.Range("C42").Value = '''https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range("C42").Style = "Normal"
$ws.Range("D42").Value = '''1.182'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '''  -3.31%  '
$ws.Range("E42").Style = "Normal"
# Row 43
$ws.Range("D43").Value = '''1.414'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '''  -4.85%  '
$ws.Range("E43").Style = "Normal"
# Row 44
$ws.Range("B44").Value = '''Frax'
$ws.Range("B44").Style = "Normal"
$ws.Range("C44").Value = '''https://coinranking.com/coin/KfWtaeV1W+frax-frax'
$ws.Range("C44").Style = "Normal"
$ws.Range("D44").Value = '''1.001'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '''  +0.08%  '
$ws.Range("E44").Style = "Normal"
# Row 45
$ws.Range("B45").Value = '''FraxShare'
$ws.Range("B45").Style = "Normal"
$ws.Range("C45").Value = '''https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range("C45").Style = "Normal"
$ws.Range("D45").Value = '''7.825'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '''  -3.57%  '
$ws.Range("E45").Style = "Normal"
# Row 46
$ws.Range("B46").Value = '''EnergySwap'
$ws.Range("B46").Style = "Normal"
$ws.Range("C46").Value = '''https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("C46").Style = "Normal"
$ws.Range("D46").Value = '''13.50'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '''  -5.31%  '
$ws.Range("E46").Style = "Normal"
# Row 47
$ws.Range("B47").Value = '''PancakeSwap'
$ws.Range("B47").Style = "Normal"
$ws.Range("C47").Value = '''https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range("C47").Style = "Normal"
$ws.Range("D47").Value = '''3.750'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '''  -3.19%  '
$ws.Range("E47").Style = "Normal"
# Row 48
$ws.Range("B48").Value = '''Decentraland'
$ws.Range("B48").Style = "Normal"
$ws.Range("C48").Value = '''https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana'
$ws.Range("C48").Style = "Normal"
$ws.Range("D48").Value = '''0.5836'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '''  -5.40%  '
$ws.Range("E48").Style = "Normal"
# Row 49
$ws.Range("B49").Value = '''Quant'
$ws.Range("B49").Style = "Normal"
$ws.Range("C49").Value = '''https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Range("C49").Style = "Normal"
$ws.Range("D49").Value = '''124.56'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '''  -3.60%  '
$ws.Range("E49").Style = "Normal"
# Row 50
$ws.Range("B50").Value = '''NEARProtocol'
$ws.Range("B50").Style = "Normal"
$ws.Range("C50").Value = '''https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range("C50").Style = "Normal"
$ws.Range("D50").Value = '''1.935'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '''  -5.82%  '
$ws.Range("E50").Style = "Normal"
# Row 51
$ws.Range("B51").Value = '''Cronos'
$ws.Range("B51").Style = "Normal"
$ws.Range("C51").Value = '''https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range("C51").Style = "Normal"
$ws.Range("D51").Value = '''0.06808'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '''  -4.48%  '
$ws.Range("E51").Style = "Normal"
